$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of resale data for 2025-01-03 09:02:32
$row = 5

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-01-03"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "09:02:32"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "00"

$ws.Cells.Item($row, 5).Value = 123656
$ws.Cells.Item($row, 6).Value = 143590
$ws.Cells.Item($row, 7).Value = 167455
$ws.Cells.Item($row, 8).Value = 157629
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 141514
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191931
$ws.Cells.Item($row, 14).Value = 114462
$ws.Cells.Item($row, 15).Value = 45212
$ws.Cells.Item($row, 16).Value = 28233
$ws.Cells.Item($row, 17).Value = 62558
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 47573
$ws.Cells.Item($row, 20).Value = -1
